$wb = $excel.ActiveWorkbook

# --- SVM sheet: add a STD (STDEV) summary row under the existing MEAN row ---
$svm = $wb.Worksheets.Item("SVM")
$svm.Range("A24").Formula = "=STDEV(A3:A12)"
$svm.Range("B24").Formula = "=STDEV(B3:B12)"

# view/selection state for SVM (best effort - scroll + active cell)
[void]$svm.Range("D19").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1

# --- GL_adaptive sheet: replace the lone "STD" label row (row 25) with a new
#     STD row at row 24 that also carries STDEV formulas across B:F ---
$gla = $wb.Worksheets.Item("GL_adaptive")
$gla.Range("A25").ClearContents()
$gla.Range("A24").Value = "STD"
$gla.Range("B24").Formula = "=STDEV(B3:B12)"
$gla.Range("C24:F24").Formula = "=STDEV(C3:C12)"

# view/selection state for GL_adaptive (keep it the active/tabSelected sheet)
[void]$gla.Range("H14").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
